# Insert a new row at position 130 (shifts existing rows 130:194 down to 131:195)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("130:130").Insert()

# Populate the newly inserted row 130 with the new weekly record
$ws.Range("A130").Value2 = 3
$ws.Range("B130").Value2 = "Femacal de La Calera"
$ws.Range("C130").Value2 = "Coquimbo"
$ws.Range("D130").Value2 = 44460
$ws.Range("E130").Value2 = 5
$ws.Range("F130").Value2 = 100112009
$ws.Range("G130").Value2 = "Acelga"
$ws.Range("H130").Value2 = "Sin especificar"
$ws.Range("I130").Value2 = "Primera"
$ws.Range("J130").Value2 = 280
$ws.Range("K130").Value2 = 1800
$ws.Range("L130").Value2 = 2000
$ws.Range("M130").Value2 = 1907
$ws.Range("N130").Value2 = "`$/docena de atados (6 kilos)"
$ws.Range("O130").Value2 = "Provincia de Quillota"
$ws.Range("P130").Value2 = 318
$ws.Range("Q130").Value2 = 6
$ws.Range("R130").Value2 = "Hortaliza"
